# Generate Report for Archive
#
# The localization status report was refreshed: the in-flight items that
# were previously marked "Ready for handoff" have moved back into
# "In Translation". Update the Status value everywhere it appears
# (the Overview rollup columns for each locale, plus the per-locale
# Status column), then let the Status column narrow to fit the new,
# shorter text - mirroring what Excel does when the column was sized to
# its contents.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: one column per locale (zh-cn -> E, de-de -> F) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E1:F1").ColumnWidth = 12.5

# --- Per-locale detail sheets: Status is column C ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C1").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C1").ColumnWidth = 12.5
